# Nalco aluminium-ingot price sheet: a new day's price row is published.
# The feed always prepends the newest date as row 2 (carrying forward the
# same description/price/circular info as the previous top row) and every
# existing row shifts down by one; the oldest row is therefore duplicated
# at the new bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Locate the extent of the existing data (row 1 = header).
$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

$firstDataRow = 2
$newLastRow = $lastRow + 1

# Work out the new top date: one day after the date currently in A2.
$topDateCell = $ws.Cells.Item($firstDataRow, 1)
$topDateText = $topDateCell.Value()
$topDate = [DateTime]::ParseExact($topDateText, "dd-MM-yyyy", $null)
$newDateText = $topDate.AddDays(1).ToString("dd-MM-yyyy")

# Shift every existing data row down by one (values + formatting),
# which also duplicates the last row into the new bottom row.
$srcRange = $ws.Range($ws.Cells.Item($firstDataRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$dstRange = $ws.Range($ws.Cells.Item($firstDataRow + 1, 1), $ws.Cells.Item($newLastRow, $lastCol))
$srcRange.Copy($dstRange)

# Build the new row 2 as a duplicate of the (now shifted) previous row 2,
# then overwrite only the date with the new, incremented date.
$newTopRowRange = $ws.Range($ws.Cells.Item($firstDataRow, 1), $ws.Cells.Item($firstDataRow, $lastCol))
$copiedTopRowRange = $ws.Range($ws.Cells.Item($firstDataRow + 1, 1), $ws.Cells.Item($firstDataRow + 1, $lastCol))
$copiedTopRowRange.Copy($newTopRowRange)

$dateCell = $ws.Cells.Item($firstDataRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = $newDateText
$ws.Cells.Item($firstDataRow + 1, 1).Copy()
$dateCell.PasteSpecial($xlPasteFormats)

# Rebuild the hyperlinks on column F (Circular Link) for every data row,
# since the link text always equals the link target in this sheet.
$ws.Hyperlinks.Delete()
for ($r = $firstDataRow; $r -le $newLastRow; $r++) {
    $linkCell = $ws.Cells.Item($r, 6)
    $url = $linkCell.Value()
    $ws.Hyperlinks.Add($linkCell, $url)
}

# Hyperlinks.Add re-styles the cell with the built-in "Hyperlink" style;
# restore the plain data-row look used throughout the sheet by copying the
# formatting from the neighbouring (unlinked) column E.
$formatSource = $ws.Range($ws.Cells.Item($firstDataRow, 5), $ws.Cells.Item($firstDataRow, 5))
$formatSource.Copy()
$linkRange = $ws.Range($ws.Cells.Item($firstDataRow, 6), $ws.Cells.Item($newLastRow, 6))
$linkRange.PasteSpecial($xlPasteFormats)

$ws.Range("A1").Select()
